$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Row 2: I2:M2 get the first 5 new peptides; N2 mirrors H2 (existing)
$ws1.Cells.Item(2,9).Value = "QLGLPADFVAK"
$ws1.Cells.Item(2,10).Value = "AVAQLVEELAR"
$ws1.Cells.Item(2,11).Value = "LTAHYCR"
$ws1.Cells.Item(2,12).Value = "QDEVSQCWNNK"
$ws1.Cells.Item(2,13).Value = "TEKPASDK"

# Rows 3-5: M column gets 3 more new peptides
$ws1.Cells.Item(3,13).Value = "QAVSDDLLK"
$ws1.Cells.Item(4,13).Value = "TLPPLNYNLVK"
$ws1.Cells.Item(5,13).Value = "SLVTEHNK"

# Rows 2-105: N column mirrors existing H column values
$ws1.Cells.Item(2,14).Value = "ALTTGVDYAQGLVALGGDDK"
$ws1.Cells.Item(3,14).Value = "AYTALLDLKPGDEFELK"
$ws1.Cells.Item(4,14).Value = "LAENAGANGAVVAENVK"
$ws1.Cells.Item(5,14).Value = "LDLNNASVR"
$ws1.Cells.Item(6,14).Value = "LGGDEYVLLSEK"
$ws1.Cells.Item(7,14).Value = "LVVGGPYSSVSDASSGLDGSQK"
$ws1.Cells.Item(8,14).Value = "QAPEVGVGDNVLYSK"
$ws1.Cells.Item(9,14).Value = "TGGDPLELFETAVK"
$ws1.Cells.Item(10,14).Value = "TNQNVGLDPETLALATPAR"
$ws1.Cells.Item(11,14).Value = "VFPGGDTEFLHPK"
$ws1.Cells.Item(12,14).Value = "VLQQEGFLSELSEEGEGVR"
$ws1.Cells.Item(13,14).Value = "YLGSTGGLLNSAETEEK"
$ws1.Cells.Item(14,14).Value = "YNSGEGGCFYSVDTLEAPWNSGR"
$ws1.Cells.Item(15,14).Value = "AAVEQLFDVR"
$ws1.Cells.Item(16,14).Value = "ALATSDLGLTPNNDGK"
$ws1.Cells.Item(17,14).Value = "ALVNNLPLYR"
$ws1.Cells.Item(18,14).Value = "ASAADTPQALAR"
$ws1.Cells.Item(19,14).Value = "DLNYYSALYEK"
$ws1.Cells.Item(20,14).Value = "LGADEYVLLSEK"
$ws1.Cells.Item(21,14).Value = "LLSQAFGLLNER"
$ws1.Cells.Item(22,14).Value = "LVADEESTSPEVSTVEEEELSALMGEK"
$ws1.Cells.Item(23,14).Value = "VPTPNVSAVDLVFESGR"
$ws1.Cells.Item(24,14).Value = "VQDYAELDGAPEER"
$ws1.Cells.Item(25,14).Value = "YVSYAVLAGDASVLQDR"
$ws1.Cells.Item(26,14).Value = "AYTALLDLK"
$ws1.Cells.Item(27,14).Value = "DAFLYYPLQYEGQECSK"
$ws1.Cells.Item(28,14).Value = "DHVPADFEK"
$ws1.Cells.Item(29,14).Value = "LQNDVQPWQVR"
$ws1.Cells.Item(30,14).Value = "SVSVSNVGTVLQLGDGLAR"
$ws1.Cells.Item(31,14).Value = "TSDTDGYAAVQLGFGDTR"
$ws1.Cells.Item(32,14).Value = "VVDADGTQLGVLNR"
$ws1.Cells.Item(33,14).Value = "AAGFALTEAEVK"
$ws1.Cells.Item(34,14).Value = "AVASGSVSAEK"
$ws1.Cells.Item(35,14).Value = "FDNADLSNANFSGAELLK"
$ws1.Cells.Item(36,14).Value = "FEQPGFFSK"
$ws1.Cells.Item(37,14).Value = "GSLPQNLGSTGGLLNSAETEEK"
$ws1.Cells.Item(38,14).Value = "LDHSQLLTDPAEAADFVAK"
$ws1.Cells.Item(39,14).Value = "LLDQDGVPVVFGGWTSASR"
$ws1.Cells.Item(40,14).Value = "LNVEYYGTETPLK"
$ws1.Cells.Item(41,14).Value = "LYLGNLPQTFESK"
$ws1.Cells.Item(42,14).Value = "QTLAMQLNEK"
$ws1.Cells.Item(43,14).Value = "WAVAEVLSNSPK"
$ws1.Cells.Item(44,14).Value = "AYTALLDLKPGDNFELK"
$ws1.Cells.Item(45,14).Value = "DMSPQALNEYK"
$ws1.Cells.Item(46,14).Value = "GFQGSNGSLFR"
$ws1.Cells.Item(47,14).Value = "LGTDAGMLAFEPSTVNLSAGDTVK"
$ws1.Cells.Item(48,14).Value = "MLTGSDLLTK"
$ws1.Cells.Item(49,14).Value = "ALQEAFQLR"
$ws1.Cells.Item(50,14).Value = "TDEEGQSLLR"
$ws1.Cells.Item(51,14).Value = "TFRPYTPGTR"
$ws1.Cells.Item(52,14).Value = "YLSYALLAGDPSVLDDR"
$ws1.Cells.Item(53,14).Value = "ACPLDVLEMVPWDGHK"
$ws1.Cells.Item(54,14).Value = "AMLPVYESK"
$ws1.Cells.Item(55,14).Value = "AQVFELPTGGAAEMNEGENLMYFAR"
$ws1.Cells.Item(56,14).Value = "DLGDADLSGSYFSVSNLQK"
$ws1.Cells.Item(57,14).Value = "EAAVADPANFDPR"
$ws1.Cells.Item(58,14).Value = "EAGFELTADEVK"
$ws1.Cells.Item(59,14).Value = "LLNYCLVTGGTGPLDELALNGQR"
$ws1.Cells.Item(60,14).Value = "SYFPYWK"
$ws1.Cells.Item(61,14).Value = "AELDYATK"
$ws1.Cells.Item(62,14).Value = "AGSTLNLDTLVK"
$ws1.Cells.Item(63,14).Value = "ANSFDDNK"
$ws1.Cells.Item(64,14).Value = "FADVVNTGK"
$ws1.Cells.Item(65,14).Value = "LLESLAPGLLK"
$ws1.Cells.Item(66,14).Value = "QWFLVDAENQTLGR"
$ws1.Cells.Item(67,14).Value = "SNQPLVNEK"
$ws1.Cells.Item(68,14).Value = "SPLANLVGWR"
$ws1.Cells.Item(69,14).Value = "SYVAAGNK"
$ws1.Cells.Item(70,14).Value = "TCSLPLDR"
$ws1.Cells.Item(71,14).Value = "AVVSADAK"
$ws1.Cells.Item(72,14).Value = "DTGVEYAQGLVALGGDDEELAK"
$ws1.Cells.Item(73,14).Value = "ETPVELEFSQLTK"
$ws1.Cells.Item(74,14).Value = "LLAQAFGLLNER"
$ws1.Cells.Item(75,14).Value = "YDSLLGQLK"
$ws1.Cells.Item(76,14).Value = "GEELELVGLR"
$ws1.Cells.Item(77,14).Value = "GFGSFEPR"
$ws1.Cells.Item(78,14).Value = "NVQATLQR"
$ws1.Cells.Item(79,14).Value = "VVTDFSEVTGR"
$ws1.Cells.Item(80,14).Value = "YGCVAGYPSGSYLGNR"
$ws1.Cells.Item(81,14).Value = "ADDEQTEENWEE"
$ws1.Cells.Item(82,14).Value = "EFTVVNVAALNELK"
$ws1.Cells.Item(83,14).Value = "LVPAGAEDSDD"
$ws1.Cells.Item(84,14).Value = "LVVGGPYASVSDASSVLDASQK"
$ws1.Cells.Item(85,14).Value = "QAGEYTTFK"
$ws1.Cells.Item(86,14).Value = "FNPGLS"
$ws1.Cells.Item(87,14).Value = "TDGLEVLSVADAAAK"
$ws1.Cells.Item(88,14).Value = "TEFDVVLEGFDAAAK"
$ws1.Cells.Item(89,14).Value = "TVTSASWR"
$ws1.Cells.Item(90,14).Value = "FADLDVSK"
$ws1.Cells.Item(91,14).Value = "FGNLVNTGR"
$ws1.Cells.Item(92,14).Value = "LFLSPVESVLR"
$ws1.Cells.Item(93,14).Value = "LGSDSGMLAFEPSSLTLQEGDTLK"
$ws1.Cells.Item(94,14).Value = "TGSTLNLDTLVK"
$ws1.Cells.Item(95,14).Value = "ACGYVSTK"
$ws1.Cells.Item(96,14).Value = "DAFLYYPLQYEAQECSNNLFYTGATPNQQSEPATK"
$ws1.Cells.Item(97,14).Value = "MLLSDLEGVTYR"
$ws1.Cells.Item(98,14).Value = "YQWDQNFYR"
$ws1.Cells.Item(99,14).Value = "ELELDDPFENLGAK"
$ws1.Cells.Item(100,14).Value = "GLVPALEAADA"
$ws1.Cells.Item(101,14).Value = "LLDQDEVPVVFGGWTSASR"
$ws1.Cells.Item(102,14).Value = "QTQAAAPVAASAEEQK"
$ws1.Cells.Item(103,14).Value = "SGEPYER"
$ws1.Cells.Item(104,14).Value = "TVYVVSDSQLEELK"
$ws1.Cells.Item(105,14).Value = "DLNALGFSDK"

# Rows 106-113: N column gets the 8 new peptides (same order as pasted above)
$ws1.Cells.Item(106,14).Value = "QLGLPADFVAK"
$ws1.Cells.Item(107,14).Value = "AVAQLVEELAR"
$ws1.Cells.Item(108,14).Value = "LTAHYCR"
$ws1.Cells.Item(109,14).Value = "QDEVSQCWNNK"
$ws1.Cells.Item(110,14).Value = "TEKPASDK"
$ws1.Cells.Item(111,14).Value = "QAVSDDLLK"
$ws1.Cells.Item(112,14).Value = "TLPPLNYNLVK"
$ws1.Cells.Item(113,14).Value = "SLVTEHNK"

# Selection / active sheet bookkeeping: select M109 on sheet1, then make the fungi sheet active
$ws1.Range("M109").Select() | Out-Null
$ws2.Activate() | Out-Null
